$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value2 = 1.070437712832501
$ws.Range("C2").Value2 = 0.1544272847382331
$ws.Range("D2").Value2 = 0.1002109442336376
$ws.Range("F2").Value2 = 3.227658738300818
$ws.Range("G2").Value2 = 0.002564428601896984
$ws.Range("J2").Value2 = 0.3563639396866449
$ws.Range("K2").Value2 = 1.146897048502979
$ws.Range("N2").Value2 = 2.978408188479008
$ws.Range("B3").Value2 = 1.023983169975196
$ws.Range("C3").Value2 = 0.1460728885522826
$ws.Range("D3").Value2 = 0.09680468056698999
$ws.Range("F3").Value2 = 3.200675288591654
$ws.Range("G3").Value2 = 0.002568925322529711
$ws.Range("J3").Value2 = 0.345809630878378
$ws.Range("K3").Value2 = 1.094593424626254
$ws.Range("N3").Value2 = 2.984968347130135
$ws.Range("B4").Value2 = 0.9961247356684453
$ws.Range("C4").Value2 = 0.1410441682637043
$ws.Range("D4").Value2 = 0.09478024708539579
$ws.Range("F4").Value2 = 3.185585212343668
$ws.Range("G4").Value2 = 0.002571831294354167
$ws.Range("J4").Value2 = 0.3395497127130227
$ws.Range("K4").Value2 = 1.063200346796833
$ws.Range("N4").Value2 = 2.989725682073413
$ws.Range("B5").Value2 = 0.9849389949196166
$ws.Range("C5").Value2 = 0.1390201696869013
$ws.Range("D5").Value2 = 0.09397213234264257
$ws.Range("F5").Value2 = 3.179806794617974
$ws.Range("G5").Value2 = 0.002573052076331739
$ws.Range("J5").Value2 = 0.3370539578197338
$ws.Range("K5").Value2 = 1.050588337664209
$ws.Range("N5").Value2 = 2.991847379353089
$ws.Range("B6").Value2 = 0.9830916744228375
$ws.Range("C6").Value2 = 0.1386856073090712
$ws.Range("D6").Value2 = 0.09383896385287471
$ws.Range("F6").Value2 = 3.178869674335687
$ws.Range("G6").Value2 = 0.002573256998956941
$ws.Range("J6").Value2 = 0.3366428677729942
$ws.Range("K6").Value2 = 1.048505035771456
$ws.Range("N6").Value2 = 2.992210730516973
$ws.Range("B7").Value2 = 0.9959732060032707
$ws.Range("C7").Value2 = 0.1410167698617926
$ws.Range("D7").Value2 = 0.09476928030223064
$ws.Range("F7").Value2 = 3.185505781734506
$ws.Range("G7").Value2 = 0.002571847610053448
$ws.Range("J7").Value2 = 0.3395158308313029
$ws.Range("K7").Value2 = 1.063029524951048
$ws.Range("N7").Value2 = 2.989753555385576
$ws.Range("B8").Value2 = 1.054281947999328
$ws.Range("C8").Value2 = 0.151525663027428
$ws.Range("D8").Value2 = 0.0990225539928673
$ws.Range("F8").Value2 = 3.218047665626486
$ws.Range("G8").Value2 = 0.002565949056387259
$ws.Range("J8").Value2 = 0.3526789143253097
$ws.Range("K8").Value2 = 1.128712666099659
$ws.Range("N8").Value2 = 2.980518529802254
$ws.Range("B9").Value2 = 1.173925428783832
$ws.Range("C9").Value2 = 0.1729414251965409
$ws.Range("D9").Value2 = 0.1078954991876486
$ws.Range("F9").Value2 = 3.29362775868168
$ws.Range("G9").Value2 = 0.002555526697383634
$ws.Range("J9").Value2 = 0.3802530769592494
$ws.Range("K9").Value2 = 1.263274244723959
$ws.Range("N9").Value2 = 2.968213202767643
$ws.Range("B10").Value2 = 1.265104391091484
$ws.Range("C10").Value2 = 0.1891798809668046
$ws.Range("D10").Value2 = 0.1147405325648521
$ws.Range("F10").Value2 = 3.356393907898962
$ws.Range("G10").Value2 = 0.002548559418574855
$ws.Range("J10").Value2 = 0.4016049869356522
$ws.Range("K10").Value2 = 1.365702891367448
$ws.Range("N10").Value2 = 2.962736460933939
$ws.Range("B11").Value2 = 1.307306310940305
$ws.Range("C11").Value2 = 0.1966793713060895
$ws.Range("D11").Value2 = 0.1179257370118592
$ws.Range("F11").Value2 = 3.386534257789862
$ws.Range("G11").Value2 = 0.002545537991415712
$ws.Range("J11").Value2 = 0.4115602826003482
$ws.Range("K11").Value2 = 1.413087844235406
$ws.Range("N11").Value2 = 2.961024500689277
$ws.Range("B12").Value2 = 1.323391867323494
$ws.Range("C12").Value2 = 0.1995355987410221
$ws.Range("D12").Value2 = 0.1191421699345909
$ws.Range("F12").Value2 = 3.398176967987666
$ws.Range("G12").Value2 = 0.002544415014223902
$ws.Range("J12").Value2 = 0.4153652295625676
$ws.Range("K12").Value2 = 1.431145672766291
$ws.Range("N12").Value2 = 2.960488763310963
$ws.Range("B13").Value2 = 1.319922894937065
$ws.Range("C13").Value2 = 0.1989197302046364
$ws.Range("D13").Value2 = 0.1188797327353655
$ws.Range("F13").Value2 = 3.395659293127949
$ws.Range("G13").Value2 = 0.002544655927770975
$ws.Range("J13").Value2 = 0.4145442011627694
$ws.Range("K13").Value2 = 1.427251508212379
$ws.Range("N13").Value2 = 2.960599131034172
$ws.Range("B14").Value2 = 1.308627580147117
$ws.Range("C14").Value2 = 0.1969140266415081
$ws.Range("D14").Value2 = 0.1180256079299227
$ws.Range("F14").Value2 = 3.38748751155012
$ws.Range("G14").Value2 = 0.002545445179821915
$ws.Range("J14").Value2 = 0.4118726130445225
$ws.Range("K14").Value2 = 1.414571180782588
$ws.Range("N14").Value2 = 2.960978167014559
$ws.Range("B15").Value2 = 1.301722510631976
$ws.Range("C15").Value2 = 0.1956876068818758
$ws.Range("D15").Value2 = 0.1175037687562934
$ws.Range("F15").Value2 = 3.382511936854314
$ws.Range("G15").Value2 = 0.00254593137337604
$ws.Range("J15").Value2 = 0.4102407661648755
$ws.Range("K15").Value2 = 1.406819001881161
$ws.Range("N15").Value2 = 2.961225007996148
$ws.Range("B16").Value2 = 1.262361008329037
$ws.Range("C16").Value2 = 0.188692048958842
$ws.Range("D16").Value2 = 0.1145338082298366
$ws.Range("F16").Value2 = 3.35445619077683
$ws.Range("G16").Value2 = 0.002548759844957499
$ws.Range("J16").Value2 = 0.4009592811915752
$ws.Range("K16").Value2 = 1.362622116073709
$ws.Range("N16").Value2 = 2.962864065229624
$ws.Range("B17").Value2 = 1.238399807045198
$ws.Range("C17").Value2 = 0.1844294301613161
$ws.Range("D17").Value2 = 0.1127301135633161
$ws.Range("F17").Value2 = 3.33765218771336
$ws.Range("G17").Value2 = 0.002550532851561983
$ws.Range("J17").Value2 = 0.3953276117049995
$ws.Range("K17").Value2 = 1.335711445907407
$ws.Range("N17").Value2 = 2.964069533392234
$ws.Range("B18").Value2 = 1.224686089453883
$ws.Range("C18").Value2 = 0.181988267751052
$ws.Range("D18").Value2 = 0.1116993930050825
$ws.Range("F18").Value2 = 3.328136357791124
$ws.Range("G18").Value2 = 0.002551566577521758
$ws.Range("J18").Value2 = 0.3921112058164056
$ws.Range("K18").Value2 = 1.320307427327407
$ws.Range("N18").Value2 = 2.964836223057532
$ws.Range("B19").Value2 = 1.220054547772122
$ws.Range("C19").Value2 = 0.1811635452265818
$ws.Range("D19").Value2 = 0.1113515622972869
$ws.Range("F19").Value2 = 3.324940091603821
$ws.Range("G19").Value2 = 0.002551918977180101
$ws.Range("J19").Value2 = 0.3910260908094045
$ws.Range("K19").Value2 = 1.315104631395258
$ws.Range("N19").Value2 = 2.965108392571793
$ws.Range("B20").Value2 = 1.240943463082829
$ws.Range("C20").Value2 = 0.1848820968058646
$ws.Range("D20").Value2 = 0.1129214247459061
$ws.Range("F20").Value2 = 3.339425534981928
$ws.Range("G20").Value2 = 0.002550342670247553
$ws.Range("J20").Value2 = 0.3959247528400169
$ws.Range("K20").Value2 = 1.338568439776765
$ws.Range("N20").Value2 = 2.963933615941627
$ws.Range("B21").Value2 = 1.311942445115733
$ws.Range("C21").Value2 = 0.197502705903446
$ws.Range("D21").Value2 = 0.1182762064667457
$ws.Range("F21").Value2 = 3.389881533811462
$ws.Range("G21").Value2 = 0.002545212783641242
$ws.Range("J21").Value2 = 0.4126563687117368
$ws.Range("K21").Value2 = 1.418292596587463
$ws.Range("N21").Value2 = 2.960863776803208
$ws.Range("B22").Value2 = 1.358954331994255
$ws.Range("C22").Value2 = 0.2058462778707337
$ws.Range("D22").Value2 = 0.1218357011071873
$ws.Range("F22").Value2 = 3.424193910318934
$ws.Range("G22").Value2 = 0.002541983458828266
$ws.Range("J22").Value2 = 0.4237961184050931
$ws.Range("K22").Value2 = 1.471062743109485
$ws.Range("N22").Value2 = 2.959513676255341
$ws.Range("B23").Value2 = 1.333807233400137
$ws.Range("C23").Value2 = 0.2013843883424897
$ws.Range("D23").Value2 = 0.1199304556045746
$ws.Range("F23").Value2 = 3.405758160000573
$ws.Range("G23").Value2 = 0.002543695760846451
$ws.Range("J23").Value2 = 0.4178318086946149
$ws.Range("K23").Value2 = 1.442837191603587
$ws.Range("N23").Value2 = 2.960174051808295
$ws.Range("B24").Value2 = 1.239793283426877
$ws.Range("C24").Value2 = 0.1846774167186993
$ws.Range("D24").Value2 = 0.1128349135011604
$ws.Range("F24").Value2 = 3.338623353034166
$ws.Range("G24").Value2 = 0.002550428606597377
$ws.Range("J24").Value2 = 0.3956547189602304
$ws.Range("K24").Value2 = 1.337276583316708
$ws.Range("N24").Value2 = 2.963994834828014
$ws.Range("B25").Value2 = 1.140986588044797
$ws.Range("C25").Value2 = 0.1670601969553616
$ws.Range("D25").Value2 = 0.1054379876790819
$ws.Range("F25").Value2 = 3.271915111714321
$ws.Range("G25").Value2 = 0.002558224484122872
$ws.Range("J25").Value2 = 0.3726029760051119
$ws.Range("K25").Value2 = 1.22624963706707
$ws.Range("N25").Value2 = 2.970918104427085
